$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep text formatting so
# numeric-looking values (e.g. "0.9900", "1.230") are not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "22.379.42"
$ws.Range("E2").Value = "  +8.90%  "
$ws.Range("D3").Value = "1.591.21"
$ws.Range("E3").Value = "  +8.01%  "
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "0.9900"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").Value = "298.76"
$ws.Range("E6").Value = "  +7.73%  "
$ws.Range("D7").Value = "0.3602"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "0.3328"
$ws.Range("E8").Value = "  +7.80%  "
$ws.Range("D9").Value = "41.02"
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("D10").Value = "1.114"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").Value = "0.06906"
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").Value = "0.9959"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "19.31"
$ws.Range("E13").Value = "  +6.27%  "
$ws.Range("D14").Value = "5.781"
$ws.Range("E14").Value = "  +5.26%  "
$ws.Range("D15").Value = "6.484"
$ws.Range("E15").Value = "  +5.08%  "
$ws.Range("D16").Value = "0.9907"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("D17").Value = "0.00001060"
$ws.Range("E17").Value = "  +3.55%  "
$ws.Range("D18").Value = "1.591.09"
$ws.Range("E18").Value = "  +8.06%  "
$ws.Range("D19").Value = "0.06570"
$ws.Range("E19").Value = "  +10.54%  "
$ws.Range("D20").Value = "76.24"
$ws.Range("E20").Value = "  +10.50%  "
$ws.Range("D21").Value = "15.79"
$ws.Range("E21").Value = "  +8.24%  "
$ws.Range("D22").Value = "5.878"
$ws.Range("E22").Value = "  +6.96%  "
$ws.Range("D23").Value = "11.45"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "22.379.49"
$ws.Range("E24").Value = "  +8.92%  "
$ws.Range("D25").Value = "2.367"
$ws.Range("E25").Value = "  +4.10%  "
$ws.Range("D26").Value = "2.504"
$ws.Range("E26").Value = "  +18.19%  "
$ws.Range("D27").Value = "148.64"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("D28").Value = "19.03"
$ws.Range("E28").Value = "  +10.74%  "
$ws.Range("D29").Value = "1.763.05"
$ws.Range("E29").Value = "  +7.98%  "
$ws.Range("D30").Value = "122.52"
$ws.Range("E30").Value = "  +7.62%  "
$ws.Range("D31").Value = "3.939"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").Value = "5.854"
$ws.Range("E32").Value = "  +18.44%  "
$ws.Range("D33").Value = "0.9174"
$ws.Range("E33").Value = "  +13.27%  "
$ws.Range("D36").Value = "11.68"
$ws.Range("E36").Value = "  +12.25%  "
$ws.Range("D37").Value = "5.075"
$ws.Range("E37").Value = "  +7.38%  "
$ws.Range("D38").Value = "1.230"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "8.336"
$ws.Range("E39").Value = "  +12.57%  "
$ws.Range("D40").Value = "0.05975"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("E41").Value = "  +6.11%  "
$ws.Range("D42").Value = "0.1966"
$ws.Range("E42").Value = "  +4.98%  "
$ws.Range("D43").Value = "0.9890"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("D44").Value = "0.5748"
$ws.Range("E44").Value = "  +9.00%  "
$ws.Range("D45").Value = "3.752"
$ws.Range("E45").Value = "  +6.68%  "
$ws.Range("D46").Value = "12.64"
$ws.Range("E46").Value = "  +3.63%  "
$ws.Range("D47").Value = "0.5577"
$ws.Range("E47").Value = "  +7.17%  "
$ws.Range("D48").Value = "123.58"
$ws.Range("E48").Value = "  +3.64%  "
$ws.Range("D49").Value = "1.928"
$ws.Range("E49").Value = "  +6.21%  "
$ws.Range("D50").Value = "0.06757"
$ws.Range("E50").Value = "  +4.69%  "
$ws.Range("D51").Value = "72.12"
$ws.Range("E51").Value = "  +7.15%  "

# Rows 34 and 35 swap coins (Stellar <-> WEMIXTOKEN) with updated data
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "1.634"
$ws.Range("E34").Value = "  +11.39%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.08087"
$ws.Range("E35").Value = "  +1.07%  "

# Restore default (Normal) style on the Price column so no stray
# cell-format/style metadata is introduced by the text-format step above.
$ws.Range("D2:E51").Style = "Normal"
